$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Locate the "Linear decay, 3rd order Runge-Kutta ordinary differential
# equation solver" row. Its whole row is removed - the "Advection &
# Diffusion" section header row and the "Uniform flow, Gaussian mass,
# concentration remote BC..." row that used to sit right after it (plus
# the now-duplicate copies of those two rows further down) collapse up
# into the gap it leaves behind.
$targetRowIndex = -1
for ($i = 1; $i -le $t.Rows.Count; $i++) {
    $cellText = $t.Cell($i, 1).Range.Text
    if ($cellText -like "*Runge-Kutta*") {
        $targetRowIndex = $i
        break
    }
}

if ($targetRowIndex -ne -1) {
    $t.Rows.Item($targetRowIndex).Delete()
} else {
    # Fallback: this row was the 15th row of the table in the source doc.
    $t.Rows.Item(15).Delete()
}

# A new, still-empty placeholder row (highlighted yellow, like the other
# "In progress" rows) is appended at the very end of the table for the
# new Zoppou-related test case that "still [does] not work".
$newRow = $t.Rows.Add()
